$d = $word.ActiveDocument

# Remove the placeholder "vnpt.SiteAddress" text that follows "Địa chỉ: "
$d.Content.Find.Execute("vnpt.SiteAddress", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
